$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.967.53"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "1.845.88"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "308.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4766"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3675"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07216"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9294"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07745"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "1.858.73"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.343"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.431"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008646"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "26.989.28"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.920"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.011"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.975"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08843"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.311"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.177"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7398"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.498"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("E35").Value = "  -4.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.111"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01959"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05265"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.976"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5242"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.024"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1514"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4736"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.010"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.605"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06070"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8888"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.59%  "
